$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data, grouped by worker: CC / EFRAIN MORALES RIVERO first (periods 2005..1909 descending),
# then PE / MOISES DAVID LINARES CORDOBA (periods 2005..1906 descending).
$data = @(
    @(16, "CC", "1047458449", "EFRAIN MORALES RIVERO", "2005", 20979),
    @(17, "CC", "1047458449", "EFRAIN MORALES RIVERO", "2004", 33125),
    @(18, "CC", "1047458449", "EFRAIN MORALES RIVERO", "2003", 33125),
    @(19, "CC", "1047458449", "EFRAIN MORALES RIVERO", "2002", 33125),
    @(20, "CC", "1047458449", "EFRAIN MORALES RIVERO", "2001", 33125),
    @(21, "CC", "1047458449", "EFRAIN MORALES RIVERO", "1912", 33125),
    @(22, "CC", "1047458449", "EFRAIN MORALES RIVERO", "1911", 33125),
    @(23, "CC", "1047458449", "EFRAIN MORALES RIVERO", "1910", 33125),
    @(24, "CC", "1047458449", "EFRAIN MORALES RIVERO", "1909", 33125),
    @(25, "PE", "949148415031996", "MOISES DAVID LINARES CORDOBA", "2005", 20979),
    @(26, "PE", "949148415031996", "MOISES DAVID LINARES CORDOBA", "2004", 33125),
    @(27, "PE", "949148415031996", "MOISES DAVID LINARES CORDOBA", "2003", 33125),
    @(28, "PE", "949148415031996", "MOISES DAVID LINARES CORDOBA", "2002", 33125),
    @(29, "PE", "949148415031996", "MOISES DAVID LINARES CORDOBA", "2001", 33125),
    @(30, "PE", "949148415031996", "MOISES DAVID LINARES CORDOBA", "1912", 33125),
    @(31, "PE", "949148415031996", "MOISES DAVID LINARES CORDOBA", "1911", 33125),
    @(32, "PE", "949148415031996", "MOISES DAVID LINARES CORDOBA", "1910", 33125),
    @(33, "PE", "949148415031996", "MOISES DAVID LINARES CORDOBA", "1909", 33125),
    @(34, "PE", "949148415031996", "MOISES DAVID LINARES CORDOBA", "1908", 33125),
    @(35, "PE", "949148415031996", "MOISES DAVID LINARES CORDOBA", "1907", 33125),
    @(36, "PE", "949148415031996", "MOISES DAVID LINARES CORDOBA", "1906", 33125)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
}
